$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("H10").NumberFormat = "@"
$ws.Range("H11").NumberFormat = "@"
$ws.Range("H12").NumberFormat = "@"
Write-Host "ok"
